$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rowsData = @(
  @("ECs","Vtn","Itga8","ECs",2,0.6666666666666666,6.597131,19.791393,0.1209543635982448,0.1209543635982448,3,1,0.8426706666666667,2.528012,0.07312747333774275,0.07312747333774275,5.559208777857333,50.03287900071599,0.008845086999114286,0.008845086999114286),
  @("ECs","Vtn","Itga8","FAPs",2,0.6666666666666666,6.597131,19.791393,0.1209543635982448,0.1209543635982448,3,1,2.843693666666667,8.531081,0.2467774671835513,0.2467774671835513,18.76021964287034,168.841976785833,0.02984881149357318,0.02984881149357317),
  @("ECs","Vtn","Itga8","M2",2,0.6666666666666666,6.597131,19.791393,0.1209543635982448,0.1209543635982448,3,1,0.430413,1.291239,0.03735150210725013,0.03735150210725013,2.839490945103,25.555418505927,0.004517827166820938,0.004517827166820938),
  @("ECs","Vtn","Itga8","sCs",2,0.6666666666666666,6.597131,19.791393,0.1209543635982448,0.1209543635982448,3,1,7.406534333333333,22.219603,0.6427435573714558,0.6427435573714558,48.86187725299767,439.756895276979,0.07774263793873637,0.07774263793873637),
  @("FAPs","Vtn","Itga8","ECs",3,1,21.05317333333333,63.15952,0.3859970617919927,0.3859970617919927,3,1,0.8426706666666667,2.528012,0.07312747333774275,0.07312747333774275,17.74089160824889,159.66802447424,0.02822698984464099,0.02822698984464099),
  @("FAPs","Vtn","Itga8","FAPs",3,1,21.05317333333333,63.15952,0.3859970617919927,0.3859970617919927,3,1,2.843693666666667,8.531081,0.2467774671835513,0.2467774671835513,59.86877567123556,538.8189810411201,0.0952553772493207,0.09525537724932069),
  @("FAPs","Vtn","Itga8","M2",3,1,21.05317333333333,63.15952,0.3859970617919927,0.3859970617919927,3,1,0.430413,1.291239,0.03735150210725013,0.03735150210725013,9.061559493919999,81.55403544528001,0.01441757006691597,0.01441757006691597),
  @("FAPs","Vtn","Itga8","sCs",3,1,21.05317333333333,63.15952,0.3859970617919927,0.3859970617919927,3,1,7.406534333333333,22.219603,0.6427435573714558,0.6427435573714558,155.9310511189511,1403.37946007056,0.2480971246311151,0.2480971246311151),
  @("sCs","Vtn","Itga8","ECs",3,1,26.89201066666667,80.676032,0.4930485746097625,0.4930485746097625,3,1,0.8426706666666667,2.528012,0.07312747333774275,0.07312747333774275,22.66110855648711,203.949977008384,0.03605539649398747,0.03605539649398747),
  @("sCs","Vtn","Itga8","FAPs",3,1,26.89201066666667,80.676032,0.4930485746097625,0.4930485746097625,3,1,2.843693666666667,8.531081,0.2467774671835513,0.2467774671835513,76.47264041673245,688.253763750592,0.1216732784406574,0.1216732784406574),
  @("sCs","Vtn","Itga8","M2",3,1,26.89201066666667,80.676032,0.4930485746097625,0.4930485746097625,3,1,0.430413,1.291239,0.03735150210725013,0.03735150210725013,11.574670987072,104.172038883648,0.01841610487351322,0.01841610487351322),
  @("sCs","Vtn","Itga8","sCs",3,1,26.89201066666667,80.676032,0.4930485746097625,0.4930485746097625,3,1,7.406534333333333,22.219603,0.6427435573714558,0.6427435573714558,199.1766002950329,1792.589402655296,0.3169037948016044,0.3169037948016044)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rowsData[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $rowVals[$j]
    }
}